$d = $word.ActiveDocument

# Merge the three runs "<id>", "p069r_2", "</id>" into a single run
# "<id>p069r_2</id>" by performing a find & replace over the exact same
# visible text. Word's Find/Replace coalesces the matched range into one
# run using the formatting of the first run in the match.
$d.Content.Find.Execute("<id>p069r_2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p069r_2</id>", 2)
